# Descriptive statistics for Galaxy J7 Duo — update the "Time to Interactive"
# (column B/E/H within each of the two summary blocks) cells so they are
# driven by a formula ( original-value * 1000 ) instead of a bare literal,
# and tighten the custom number format from 6 to 3 decimal places so the
# on-screen text is unchanged even though the stored number is 1000x larger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> original literal value (kept exactly as typed originally)
# so the rebuilt formula text matches "=<value>*1000".
$formulas = [ordered]@{
    "B4"  = "0.008447*1000"
    "E4"  = "0.0088*1000"
    "H4"  = "0.007798*1000"
    "B5"  = "0.003791*1000"
    "E5"  = "0.005402*1000"
    "H5"  = "0.003233*1000"
    "B6"  = "0.01679*1000"
    "E6"  = "0.023197*1000"
    "H6"  = "0.014613*1000"
    "B7"  = "0.010273*1000"
    "E7"  = "0.010927*1000"
    "H7"  = "0.0104*1000"
    "B8"  = "0.007795*1000"
    "E8"  = "0.00703*1000"
    "H8"  = "0.006806*1000"
    "B9"  = "0.005421*1000"
    "E9"  = "0.00499*1000"
    "H9"  = "0.005683*1000"
    "B10" = "0.00396*1000"
    "E10" = "0.003832*1000"
    "H10" = "0.002169*1000"
    "B15" = "0.00818*1000"
    "E15" = "0.008895*1000"
    "H15" = "0.0079*1000"
    "B16" = "0.003482*1000"
    "E16" = "0.004515*1000"
    "H16" = "0.004532*1000"
    "B17" = "0.0157*1000"
    "E17" = "0.01898*1000"
    "H17" = "0.02115*1000"
    "B18" = "0.0101*1000"
    "E18" = "0.01074*1000"
    "H18" = "0.01023*1000"
    "B19" = "0.00712*1000"
    "E19" = "0.00781*1000"
    "H19" = "0.0064*1000"
    "B20" = "0.00553*1000"
    "E20" = "0.00557*1000"
    "H20" = "0.00517*1000"
    "B21" = "0.00425*1000"
    "E21" = "0.00332*1000"
    "H21" = "0.00262*1000"
}

foreach ($ref in $formulas.Keys) {
    $cell = $ws.Range($ref)

    # Values moved from a bare literal (e.g. 0.008447) to a formula that
    # multiplies the original reading by 1000.
    $cell.Formula = "=" + $formulas[$ref]

    # The stored magnitude grew 1000x, so drop three decimal places from the
    # custom display format (0.000000 -> 0.000) to keep the on-screen text
    # identical to what it was before the edit.
    $cell.NumberFormat = "0.000"
}

# Cursor/selection left on E27 by the author when the file was last saved.
$ws.Range("E27").Select()
